# Insert a new product row at row 7 ("Pepsi" / "dietética" / "black"),
# shifting the existing rows 7-38 down to 8-39, preserving all data and
# formatting. This mirrors an "insert row above" operation followed by
# filling in the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data block (rows 7-38) down by one row, as a single
# block copy/paste so that values, number formats and styles all move
# together (avoids Excel's "insert row" auto-formatting side effects).
$srcBlock = $ws.Range("A7:O38")
$dstBlock = $ws.Range("A8:O39")
$srcBlock.Copy()
$dstBlock.PasteSpecial(-4104)   # xlPasteAll

# The newly created row 39 doesn't automatically inherit column A's
# "0" number format (used to avoid scientific notation for long
# barcodes) since it previously fell outside the worksheet's used
# range, so make sure it matches the rest of the column explicitly.
$ws.Range("A39").NumberFormat = "0"

$ws.Application.CutCopyMode = $false

# Fill in the new row 7 with the new product's data.
$ws.Range("A7").Value = 7791813420057
$ws.Range("B7").Value = "Gaseosa"
$ws.Range("C7").Value = "dietética"
$ws.Range("D7").Value = "black"
$ws.Range("E7").Value = "Pepsi"
$ws.Range("F7").Value = 500
$ws.Range("G7").Value = "ml."
$ws.Range("H7").Value = "Botella"
$ws.Range("I7").Value = "Gaseosas"
$ws.Range("J7").Value = "Argentina"
$ws.Range("K7").Value = 6
$ws.Range("L7").Value = $false
$ws.Range("M7").Value = $true
$ws.Range("N7").Value = "C:\VentaSoft\Imágenes de artículos\7791813420057.png"
$ws.Range("O7").Value = $true
